$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Main" table by one row (A1:H66 -> A1:H67)
$lo = $ws.ListObjects.Item("Main")
$lo.ListRows.Add() | Out-Null
$lo.Resize($ws.Range("A1:H67"))

# Copy number formats/styles from row 66 down to the new row 67
$ws.Range("D66:H66").Copy()
$ws.Range("D67").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new part row - order matters for shared-string insertion order
$ws.Range("H67").Value = "Allegro (msalamon_pl)"
$ws.Range("G67").Value = "https://allegro.pl/oferta/wskaznik-naladowania-akumulatora-li-ion-uniwersalny-zielony-32-x-20-mm-13040945266"
$ws.Range("C67").Value = "Ordered"
$ws.Range("A67").Value = "Lithium battery charge indicator (1S-8S, green)"
$ws.Range("B67").Value = 1
$ws.Range("D67").Value = 9.6
$ws.Range("E67").Formula = "=PRODUCT(B67*D67)"
$ws.Range("F67").Value = "2/1/2024"

# Turn the LINK cell into a real hyperlink
$ws.Hyperlinks.Add($ws.Range("G67"), "https://allegro.pl/oferta/wskaznik-naladowania-akumulatora-li-ion-uniwersalny-zielony-32-x-20-mm-13040945266") | Out-Null
$ws.Range("G67").Style = "Hyperlink"

# Update the view selection to the new last cell
$ws.Range("H68").Select() | Out-Null
